# TODO.xlsx update:
#  - Finished "Scene transitions" work and added narration status
#  - Bold the header row
#  - Insert new "Transitions" rows under the Intro/Audio block
#  - Mark narration-related rows as Done
#  - Insert a new "Trigger Stay Tuned panel" gameplay row near the end
#  - Fill in a couple of previously-blank Status cells
#  - Update selection / page orientation to match the author's last save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make the header row (row 1) bold ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Intro/Audio block: narration & friends are now Done ---
$ws.Range("D5").Value = "Done"
$ws.Range("D6").Value = "Done"
$ws.Range("D7").Value = "Done"

# --- Insert 3 new rows (9:11) for the new "Transitions" section, pushing the ---
# --- "Main" section (old row 10 onward) down by three rows.                 ---
$ws.Rows("9:11").Insert()

$ws.Range("B9").Value = "Transitions"
$ws.Range("C9").Value = "Menu to Instructions"
$ws.Range("D9").Value = "Done"

$ws.Range("B10").Value = "Transitions"
$ws.Range("C10").Value = "Menu to Intro"
$ws.Range("D10").Value = "Done"

$ws.Range("B11").Value = "Transitions"
$ws.Range("C11").Value = "Intro to Chapter 1"
$ws.Range("D11").Value = "Done!!"

# --- The "Door handle on inside of top door" row (now row 22) is Done too ---
$ws.Range("D22").Value = "Done"

# --- Insert 1 new row (43) for the new "Trigger Stay Tuned panel" gameplay ---
# --- item, pushing the final "Credits scene" row down to row 44.          ---
$ws.Rows("43:43").Insert()

# --- "Animation of middle cask opening" (now row 42) is Done ---
$ws.Range("D42").Value = "Done"

$ws.Range("B43").Value = "Gameplay"
$ws.Range("C43").Value = "Trigger Stay Tuned panel"
$ws.Range("D43").Value = "Done"

# --- Match the author's final selection / scroll position ---
$ws.Range("D11").Select()

# --- Touch page setup (orientation) so it gets serialized like the source ---
$ws.PageSetup.Orientation = 1
